$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-VertexCoordText($ShapeName, $NewText) {
    $shape = $s.Shapes.Item($ShapeName)
    $paragraph = $shape.TextFrame.TextRange.Paragraphs(2)
    $paragraph.Runs(1).Text = $NewText
}

Set-VertexCoordText "TextBox 10" "( -583, 81)"
Set-VertexCoordText "TextBox 54" "( -486, 81)"
Set-VertexCoordText "TextBox 55" "( -426,81)"
Set-VertexCoordText "TextBox 56" "( -306, 81)"
Set-VertexCoordText "TextBox 57" "( -364, 81)"
Set-VertexCoordText "TextBox 58" "( -248, 81)"
Set-VertexCoordText "TextBox 59" "( -130, 81)"
Set-VertexCoordText "TextBox 60" "( -189, 81)"
Set-VertexCoordText "TextBox 61" "(45, 81)"
Set-VertexCoordText "TextBox 62" "( 173, 81)"
Set-VertexCoordText "TextBox 63" "( 108, 81)"
Set-VertexCoordText "TextBox 64" "( 230, 81)"
Set-VertexCoordText "TextBox 65" "( 350, 81)"
Set-VertexCoordText "TextBox 66" "( 290, 81)"
Set-VertexCoordText "TextBox 67" "( 406, 81)"
Set-VertexCoordText "TextBox 68" "( 466, 81)"
Set-VertexCoordText "TextBox 69" "( 526, 81)"
Set-VertexCoordText "TextBox 108" "( -583, 87)"
Set-VertexCoordText "TextBox 152" "( -486, 87)"
Set-VertexCoordText "TextBox 153" "( -426,87)"
Set-VertexCoordText "TextBox 154" "( -306, 87)"
Set-VertexCoordText "TextBox 155" "( -364, 87)"
Set-VertexCoordText "TextBox 156" "( -248, 87)"
Set-VertexCoordText "TextBox 157" "( -130, 87)"
Set-VertexCoordText "TextBox 158" "( -189, 87)"
Set-VertexCoordText "TextBox 159" "(45, 87)"
Set-VertexCoordText "TextBox 160" "( 173, 87)"
Set-VertexCoordText "TextBox 161" "( 108, 87)"
Set-VertexCoordText "TextBox 162" "( 230, 87)"
Set-VertexCoordText "TextBox 163" "( 350, 87)"
Set-VertexCoordText "TextBox 164" "( 290, 87)"
Set-VertexCoordText "TextBox 165" "( 406, 87)"
Set-VertexCoordText "TextBox 166" "( 466, 87)"
Set-VertexCoordText "TextBox 167" "( 526, 87)"
Set-VertexCoordText "TextBox 186" "( 583, 81)"
Set-VertexCoordText "TextBox 190" "( 583, 87)"
